$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: "Base" -> "Base (0%)", values updated ---
$ws.Range("A2").Value = "Base (0%)"
$ws.Range("B2").Value = 43956
$ws.Range("C2").Value = 1489
$ws.Range("D2").Value = 25
$ws.Range("E2").Value = 43962
$ws.Range("F2").Value = 1567
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = 3018
$ws.Range("I2").Value = 463

# --- Row 3: "20% reduction" label unchanged, values updated ---
$ws.Range("A3").Value = "20% reduction"
$ws.Range("B3").Value = 43958
$ws.Range("C3").Value = 1451
$ws.Range("D3").Value = 25
$ws.Range("E3").Value = 43963
$ws.Range("F3").Value = 1534
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 3010
$ws.Range("I3").Value = 469

# --- Row 4: now "40% reduction" (new scenario row), values updated ---
$ws.Range("A4").Value = "40% reduction"
$ws.Range("B4").Value = 43958
$ws.Range("C4").Value = 1459
$ws.Range("D4").Value = 24
$ws.Range("E4").Value = 43963
$ws.Range("F4").Value = 1542
$ws.Range("G4").Value = 32
$ws.Range("H4").Value = 3011
$ws.Range("I4").Value = 468

# --- Row 5: new row "50% reduction" ---
$ws.Range("A5").Value = "50% reduction"
$ws.Range("B5").Value = 43966
$ws.Range("C5").Value = 1259
$ws.Range("D5").Value = 23
$ws.Range("E5").Value = 43971
$ws.Range("F5").Value = 1358
$ws.Range("G5").Value = 29
$ws.Range("H5").Value = 2993
$ws.Range("I5").Value = 503

# --- Row 6: new row "60% reduction" ---
$ws.Range("A6").Value = "60% reduction"
$ws.Range("B6").Value = 43974
$ws.Range("C6").Value = 1110
$ws.Range("D6").Value = 22
$ws.Range("E6").Value = 43979
$ws.Range("F6").Value = 1211
$ws.Range("G6").Value = 29
$ws.Range("H6").Value = 2962
$ws.Range("I6").Value = 522

# --- Row 7: new row "73% reduction" ---
$ws.Range("A7").Value = "73% reduction"
$ws.Range("B7").Value = 43988
$ws.Range("C7").Value = 835
$ws.Range("D7").Value = 18
$ws.Range("E7").Value = 43993
$ws.Range("F7").Value = 932
$ws.Range("G7").Value = 25
$ws.Range("H7").Value = 2879
$ws.Range("I7").Value = 597

# Apply the existing date cell style (same as B2/E2 already have, which maps
# to the workbook's date number format) to the newly-added date cells in
# rows 4-7 so they share the same style index rather than creating new ones.
$dateCells = @("B4", "E4", "B5", "E5", "B6", "E6", "B7", "E7")
foreach ($addr in $dateCells) {
    $ws.Range("B2").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
